$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "LeaveEmpty" column to the menu table (extends table ref
# A1:F5 -> A1:G5 and inserts a new ListColumn at the end).
$lo = $ws.ListObjects.Item("Table3")
$lo.ListColumns.Add() | Out-Null
$ws.Range("G1").Value = "LeaveEmpty"

# House Salad (row 4) previously had no Allergens entry - fill it in to
# match the other rows ("No known priority allergens"), with the same
# wrapped-text styling used by the other Allergens cells.
$ws.Range("C4").Value = "No known priority allergens"
$ws.Range("C4").WrapText = $true

# Update the active selection to reflect where the edit left off.
$ws.Range("C4").Select() | Out-Null
